# Applies the "update LCOM composition calculation" edit:
#  - Calculation sheet: add fom-PV helper block (E12:F14, L15:P21, B10/B21 blanks)
#  - Calculation!B8 becomes a formula (=E13+F14) instead of a literal
#  - Calculation!C8 / D8 become plain (non-shared) formulas
#  - N2:N8 / O2:O8 formulas re-entered as shared formulas
#  - Selection / active-tab bookkeeping moves from LCOM_composition to Calculation
#  - Recalculated cached values ripple into LCOM_composition (handled by recalc)

$wb = $excel.ActiveWorkbook
$calc = $wb.Worksheets.Item("Calculation")
$lcom = $wb.Worksheets.Item("LCOM_composition")

# ---------------------------------------------------------------------------
# 1. New shared strings, inserted in the order the diff expects them to land
#    in sharedStrings.xml (unit=46, fom PV=47, fom rest after variable cost=48)
# ---------------------------------------------------------------------------
$calc.Range("M15").Value = "unit"
$calc.Range("M16").Value = "unit"
$calc.Range("M17").Value = "unit"

$calc.Range("E12").Value = "fom PV"
$calc.Range("F12").Value = "fom rest after variable cost"

# ---------------------------------------------------------------------------
# 2. New small lookup table (rows 15-17, cols L:N) that labels the fom-PV
#    inputs, plus the per-unit fom figures themselves (col O) and the PV
#    total (P15)
# ---------------------------------------------------------------------------
$calc.Range("L15").Value = "electrolyzer"
$calc.Range("N15").Value = "fom_cost"
$calc.Range("O15").Value = 4.34
$calc.Range("P15").Formula = '=O15*8254*44'

$calc.Range("L16").Value = "ch3oh_reactor"
$calc.Range("N16").Value = "fom_cost"
$calc.Range("O16").Value = 4.45

$calc.Range("L17").Value = "steam_plant"
$calc.Range("N17").Value = "fom_cost"
$calc.Range("O17").Value = 0.1192922374429224

# ---------------------------------------------------------------------------
# 3. Utilisation-style factors (L19:L21) feeding O19:O21
# ---------------------------------------------------------------------------
$calc.Range("L19").Value = 0.82802811702252499
$calc.Range("O19").Formula = '=L19*52*O15'

$calc.Range("L20").Value = 0.468957494005449
$calc.Range("O20").Formula = '=L20*52*O16'

$calc.Range("L21").Value = 0.0175757509509804
$calc.Range("O21").Formula = '=L21*100*O17'

# ---------------------------------------------------------------------------
# 4. fom PV / fom rest after variable cost block (E13:F14)
# ---------------------------------------------------------------------------
$calc.Range("E13").Formula = '=304*11300'
$calc.Range("F13").Value = 1921371.8921226684
$calc.Range("F13").NumberFormat = "#,##0.00"
$calc.Range("F13").Font.Name = "Arial Unicode MS"
$calc.Range("F13").Font.Size = 10
$calc.Range("F13").VerticalAlignment = -4108
$calc.Range("F14").Formula = '=F13-SUM(B2:B7)'

# ---------------------------------------------------------------------------
# 5. Calculation!B8 becomes a formula; C8/D8 become plain (non-shared)
#    formulas referencing the new B8
# ---------------------------------------------------------------------------
$calc.Range("B8").Formula = '=E13+F14'
$calc.Range("C8").Formula = '=(B8*$G$4)/($G$2*$G$4)'
$calc.Range("D8").Formula = '=(B8*$G$4)/($H$2*$G$4)'

# ---------------------------------------------------------------------------
# 6. Re-enter N2:N8 / O2:O8 as a single fill so they come back as shared
#    formulas (matches the diff's t="shared" restructuring)
# ---------------------------------------------------------------------------
$calc.Range("N2:N8").Formula = '=M2/($G$2*$G$4)'
$calc.Range("O2:O8").Formula = '=M2/($H$2*$G$4)'

# ---------------------------------------------------------------------------
# 7. Blank-but-styled cells added by the edit
# ---------------------------------------------------------------------------
$calc.Range("B10").NumberFormat = "#,##0.00"
$calc.Range("B10").Font.Name = "Arial Unicode MS"
$calc.Range("B10").Font.Size = 10
$calc.Range("B10").VerticalAlignment = -4108

$calc.Range("B21").NumberFormat = "#,##0.00"

$calc.Range("L20").Font.Name = "Arial Unicode MS"
$calc.Range("L20").Font.Size = 10
$calc.Range("L20").VerticalAlignment = -4108

$calc.Range("L21").Font.Name = "Arial Unicode MS"
$calc.Range("L21").Font.Size = 10
$calc.Range("L21").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 8. Column widths touched by the edit (ColumnWidth is quantised to 1/6-char
#    steps by this host, so these inputs are chosen to land on the closest
#    reproducible width to the authored one)
# ---------------------------------------------------------------------------
$calc.Columns.Item(6).ColumnWidth = 11.59
$calc.Columns.Item(12).ColumnWidth = 14.42
$calc.Columns.Item(16).ColumnWidth = 11.09

# ---------------------------------------------------------------------------
# 9. Selection / active tab bookkeeping: the edit moves the active sheet
#    from LCOM_composition to Calculation, with fresh selections on both.
# ---------------------------------------------------------------------------
$lcom.Range("G20").Select()
$calc.Activate()
$calc.Range("F13").Select()

$wb.Save()
